$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update A8 (Dashboard Page with trailing space - new shared string)
$ws.Range("A8").Value = "Dashboard Page "

# 2. Update A9 (Dashboard Page -> Product listing page - new shared string)
$ws.Range("A9").Value = "Dashboard Page -> Product listing page"

# 3. Fill row 12 - new Product details page test case
$ws.Range("A12").Value = "product details page"
$ws.Range("B12").Value = "TC1"
$ws.Range("C12").Value = "verify page title on product details page"
$ws.Range("D12").Value = "Medium"
$ws.Range("E12").Value = "Sanity"

# 4. Fill row 13 - another Product details page test case
$ws.Range("A13").Value = "product details page"
$ws.Range("B13").Value = "TC2"
$ws.Range("C13").Value = "verify user redirects to correct details page"
$ws.Range("D13").Value = "Medium"
$ws.Range("E13").Value = "Sanity"

# 5. Update E4 (Sanity, Regression - new shared string)
$ws.Range("E4").Value = "Sanity, Regression"

# 6. Resize columns: A wider to fit new text, split D:E width so E can be wider
$ws.Columns.Item(1).ColumnWidth = 36
$ws.Columns.Item(5).ColumnWidth = 19.166666666666668

# 7. Update selection to C16
$ws.Range("C16").Select() | Out-Null

Write-Host "done"
